# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E22) listed the last 7 billing periods in
# descending order (2409 .. 2403). The database refresh re-sorts them in
# ascending order (2403 .. 2409), which moves the "Valor Mora" amount that
# belongs to period 2409 (32933) down to the row that now holds 2409, while
# every other period keeps its 52000 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the period labels in column E (rows 16-22) from descending to
# ascending order. Row 19 (2406) is the midpoint and is unchanged.
$ws.Range("E16").Value = "2403"
$ws.Range("E17").Value = "2404"
$ws.Range("E18").Value = "2405"
$ws.Range("E19").Value = "2406"
$ws.Range("E20").Value = "2407"
$ws.Range("E21").Value = "2408"
$ws.Range("E22").Value = "2409"

# The "Valor Mora" figure that was tied to period 2409 travels with it: it
# used to sit in row 16 (32933) and now belongs in row 22, where 2409 now
# lives. Row 16 picks up the standard 52000 value instead.
$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 32933
